$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace two names in the signature list (rows 2 and 3, column A)
$ws.Range("A2").Value = "DEBAGH OUSSAMA"
$ws.Range("A3").Value = "SBAIHI SID ALI"

# Update the active selection shown when the file was last saved
$ws.Range("A8:E8").Select() | Out-Null
